# Edit: split the former CASH_FLOWS sheet into two statements.
# 1) The existing "CASH_FLOWS" sheet (sheetId=2) is renamed to
#    "CONSOLIDATED_BALANCE_SHEETS" and its data is replaced with the
#    consolidated balance sheet figures (from page 5 of the 10-Q).
# 2) A brand new "CASH_FLOWS" sheet is inserted right after it (sheetId=3)
#    and repopulated with the original cash-flow-statement figures
#    (from page 9 of the 10-Q), including the matching header comment.

$wb = $excel.ActiveWorkbook

$balanceSheetData = (
        @('Description', 'Year Ended December 31, 2024', 'Year Ended December 31, 2025'),
        @('Cash and cash equivalents', '23466', '23264'),
        @('Marketable securities', '72191', '72064'),
        @('Total cash, cash equivalents, and marketable securities', '95657', '95328'),
        @('Accounts receivable, net', '52340', '51000'),
        @('Other current assets', '15714', '15724'),
        @('Total current assets', '163711', '162052'),
        @('Non-marketable securities', '37982', '51029'),
        @('Deferred income taxes', '17180', '18386'),
        @('Property and equipment, net', '171036', '185062'),
        @('Operating lease assets', '13588', '13722'),
        @('Goodwill', '31885', '32173'),
        @('Other non-current assets', '14874', '12950'),
        @('Total assets', '450256', '475374'),
        @('Accounts payable', '7987', '8497'),
        @('Accrued compensation and benefits', '15069', '9984'),
        @('Accrued expenses and other current liabilities', '51228', '58300'),
        @('Accrued revenue share', '9802', '9965'),
        @('Deferred revenue', '5036', '4908'),
        @('Total current liabilities', '89122', '91654'),
        @('Long-term debt', '10883', '10886'),
        @('Income taxes payable, non-current', '8782', '9773'),
        @('Operating lease liabilities', '11691', '11678'),
        @('Other long-term liabilities', '4694', '6116'),
        @('Total liabilities', '125172', '130107'),
        @('Accumulated other comprehensive income (loss)', '-4800', '-4086'),
        @('Retained earnings', '245084', '262628'),
        @('Total stockholders’ equity', '325084', '345267'),
        @('Total liabilities and stockholders’ equity', '450256', '475374')
    )

$cashFlowsData = (
        @('Description', 'Year Ended December 31, 2024', 'Year Ended December 31, 2025'),
        @('Net income', '23662', '34540'),
        @('Depreciation of property and equipment', '3413', '4487'),
        @('Stock-based compensation expense', '5264', '5516'),
        @('Deferred income taxes', '419', '-1152'),
        @('Loss (gain) on debt and equity securities, net', '-1781', '-9960'),
        @('Other', '334', '481'),
        @('Accounts receivable, net', '3167', '1638'),
        @('Income taxes, net', '3011', '7197'),
        @('Other assets', '-1000', '-1288'),
        @('Accounts payable', '-2124', '-880'),
        @('Accrued expenses and other liabilities', '-5054', '-5045'),
        @('Accrued revenue share', '-322', '116'),
        @('Deferred revenue', '-141', '500'),
        @('Net cash provided by operating activities', '28848', '36150'),
        @('Purchases of property and equipment', '-12012', '-17197'),
        @('Purchases of marketable securities', '-20684', '-18453'),
        @('Maturities and sales of marketable securities', '24985', '20345'),
        @('Purchases of non-marketable securities', '-1206', '-958'),
        @('Maturities and sales of non-marketable securities', '313', '259'),
        @('Acquisitions, net of cash acquired, and purchases of intangible assets', '-61', '-340'),
        @('Other investing activities', '101', '150'),
        @('Net cash used in investing activities', '-8564', '-16194'),
        @('Net payments related to stock-based award activities', '-2929', '-3110'),
        @('Repurchases of stock', '-15696', '-15068'),
        @('Dividend payments', '0', '-2434'),
        @('Proceeds from issuance of debt, net of costs', '1982', '4532'),
        @('Repayments of debt', '-3079', '-4521'),
        @('Proceeds from sale of interest in consolidated entities, net', '8', '400'),
        @('Net cash used in financing activities', '-19714', '-20201'),
        @('Effect of exchange rate changes on cash and cash equivalents', '-125', '43'),
        @('Net increase (decrease) in cash and cash equivalents', '445', '-202'),
        @('Cash and cash equivalents at beginning of period', '24048', '23466'),
        @('Cash and cash equivalents at end of period', '24493', '23264')
    )

# Writes a value into a cell. Numeric-looking strings are written with a
# leading apostrophe (same as typing '23466 into Excel) so they land as
# text, not numbers -- matching the source PDF-extraction's inlineStr
# cells -- and then the quote-prefix formatting flag is cleared so the
# cell's style stays the plain default (no stray "number stored as text"
# styling left behind).
function Set-CellValue {
    param($cell, $val)
    if ($val -match '^-?\d+(\.\d+)?$') {
        $cell.Value = "'" + $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

function Write-GridData {
    param($sheet, $grid)
    for ($r = 0; $r -lt $grid.Count; $r++) {
        $row = $grid[$r]
        for ($c = 0; $c -lt $row.Count; $c++) {
            Set-CellValue $sheet.Cells.Item($r + 1, $c + 1) $row[$c]
        }
    }
}

# --- Step 1: repurpose the current CASH_FLOWS sheet as the balance sheet ---
$balanceSheet = $wb.Worksheets.Item("CASH_FLOWS")
$balanceSheet.Name = "CONSOLIDATED_BALANCE_SHEETS"

# Wipe the old cash-flow-statement content before writing the new data.
$balanceSheet.Cells.Clear()

Write-GridData $balanceSheet $balanceSheetData

# Match the header styling used on the other statement sheets (bold,
# centered, top-aligned, thin border all around).
$incomeSheet = $wb.Worksheets.Item("INCOME")
$incomeSheet.Range("A1:C1").Copy()
$balanceSheet.Range("A1:C1").PasteSpecial(-4122)

# Update the page-reference comment for the balance sheet.
$balanceSheet.Range("A1").AddComment("From page 5 of goog-10-q-q1-2025.pdf")

$balanceSheet.Columns.Item(1).AutoFit()
$balanceSheet.Columns.Item(2).AutoFit()
$balanceSheet.Columns.Item(3).AutoFit()

# --- Step 2: insert a new CASH_FLOWS sheet right after the balance sheet ---
$cashFlowSheet = $wb.Worksheets.Add($null, $balanceSheet)
$cashFlowSheet.Name = "CASH_FLOWS"

Write-GridData $cashFlowSheet $cashFlowsData

$incomeSheet.Range("A1:C1").Copy()
$cashFlowSheet.Range("A1:C1").PasteSpecial(-4122)

$cashFlowSheet.Range("A1").AddComment("From page 9 of goog-10-q-q1-2025.pdf")

$cashFlowSheet.Columns.Item(1).AutoFit()
$cashFlowSheet.Columns.Item(2).AutoFit()
$cashFlowSheet.Columns.Item(3).AutoFit()

$wb.Worksheets.Item("INCOME").Select()

Write-Output "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
